# Update workbook "market_health_data.xlsx":
#  1. Metadata!A2 timestamp bump (10:44 -> 10:45)
#  2. Industry Analysis!F2:F76 - refreshed "1 Year" return figures
#  3. Stock List - a new stock (CAPTRU-RE1) appears at the top of the list,
#     pushing every other row down by one and dropping the prior last row
#     (TRAVELFOOD) off the bottom.

$wb = $excel.ActiveWorkbook

# 1. Metadata timestamp -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 10:45 AM"

# 2. Industry Analysis "1 Year" column (F) -------------------------------
$industry = $wb.Worksheets.Item("Industry Analysis")
$newFValues = @(
    18.476, -7.7404, 30.7972, -50.2266, 61.9649, -9.1713, -3.556, 38.3509,
    -6.2497, 52.6723, -6.932, 17.5662, -35.5106, 0.6286, -3.1514, -20.6354,
    -0.0175, -26.9255, 44.703, 10.0506, 84.6016, -54.4868, -12.8122,
    -9.182700000000001, 5.9529, -33.2998, -20.4441, -17.1514, 24.527,
    57.6193, -1.527, -5.2378, 27.4054, 6.7961, -5.6683, 1.4178, -22.4272,
    12.3741, -5.138, -0.1825, 23.2483, 14.456, -11.1739, 27.112, -5.6252,
    -36.5148, -27.8397, -25.4424, -49.1173, -51.065, -35.4517, -11.9879,
    -3.0992, -15.3441, -25.937, -29.1486, -6.4093, -23.3046, -11.2657,
    -9.777699999999999, -16.0561, -9.932499999999999, 51.8767, -43.5191,
    13.7315, 12.6111, 31.7532, -19.9577, -12.9642, 13.2432, 2.8232, -9.179,
    -14.2931, 28.3699, 45.5868
)
for ($i = 0; $i -lt $newFValues.Length; $i++) {
    $industry.Cells.Item($i + 2, 6).Value = $newFValues[$i]
}

# 3. Stock List - insert the new top row, then drop the old bottom row --
$stocks = $wb.Worksheets.Item("Stock List")
$stocks.Rows.Item(2).Insert()
$stocks.Rows.Item(2).ClearFormats()
$stocks.Range("A2").Value = [char]0x1F4CB
$stocks.Range("B2").Value = "CAPTRU-RE1"
$stocks.Range("C2").Value = "CAPTRU-RE1"
$stocks.Range("D2").Value = 5.67
$stocks.Range("E2").Value = -11.9565
$stocks.Range("F2").Value = "N/A"
$stocks.Range("G2").Value = "N/A"
$stocks.Range("H2").Value = 0
$stocks.Rows.Item(77).Delete()
